# Add annotations.de / annotations.en columns for all projects
# (surveys.xlsx, sheet "surveys")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells — these append two new shared strings and extend the
# sheet's used range from A1:N3 to A1:P3 automatically.
$ws.Range("O1").Value = "annotations.de"
$ws.Range("P1").Value = "annotations.en"

# Give the new "annotations.de" column a bit more breathing room.
$ws.Columns.Item(15).ColumnWidth = 17.3

# Scroll the view over so column C is left-most and select the first data
# cell in the newly added "annotations.en" column.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("P2").Select()
